$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.211.72'
$ws.Range('E2').Value = '  -0.25%  '
$ws.Range('D3').Value = '1.682.42'
$ws.Range('E3').Value = '  +0.35%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '216.11'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.51%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5257'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.52%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2696'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.46%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06365'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.61%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '21.51'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.71%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07623'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.58%  '
$ws.Range('D12').Value = '1.719.12'
$ws.Range('E12').Value = '  +2.38%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.515'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.09%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5759'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.11%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.000008330'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.94%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '66.09'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.34%  '
$ws.Range('D17').Value = '26.244.59'
$ws.Range('E17').Value = '  -0.26%  '
$ws.Range('E18').Value = '  +0.02%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.864'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.86%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '189.33'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.22%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.230'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.82%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.008'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.03%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '148.84'
$ws.Range('D24').Style = 'Normal'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '7.788'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.17%  '
$ws.Range('E26').Value = '  -0.64%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.06305'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.88%  '
$ws.Range('E29').Value = '  +0.85%  '
$ws.Range('E30').Value = '  -0.11%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.563'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.58%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.569'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.50%  '
$ws.Range('E33').Value = '  +1.79%  '
$ws.Range('E34').Value = '  -0.65%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6116'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.20%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.420'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.67%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.756'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.06%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.185'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.51%  '
$ws.Range('E39').Value = '  -0.19%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.8918'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.12%  '
$ws.Range('D41').Value = '1.099.41'
$ws.Range('E41').Value = '  -1.30%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '100.47'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.00%  '
$ws.Range('E44').Value = '  +0.33%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00000000110'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.34%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '57.30'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.73%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.006'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.27%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.042'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.77%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05275'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.30%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.4286'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.10%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.016'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.71%  '
